$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F35").Value = "https://www.youtube.com/watch?v=sfHG9Vg-ELE"
$ws.Range("B61").Value = "Adicto"
$ws.Range("C61").Value = "Carlitos Rossy, Dalex"
$ws.Range("D61").Value = "Sad Valentín"
$ws.Range("E61").Value = "'2023-05-25"
$ws.Range("F61").Value = "https://www.youtube.com/watch?v=URvNDDbLIZ8"
$ws.Range("F78").Value = "https://www.youtube.com/watch?v=EY-nAOFpenI"
$ws.Range("F88").Value = "https://www.youtube.com/watch?v=6HgF6T02ZBY"
$ws.Range("F114").Value = "https://www.youtube.com/watch?v=D6Ju9CyOB-I"
$ws.Range("F141").Value = "https://www.youtube.com/watch?v=WENoupAz5C0"
$ws.Range("B161").Value = "SUERTE"
$ws.Range("C161").Value = "Renn"
$ws.Range("D161").Value = "SUERTE"
$ws.Range("E161").Value = "'2023-03-23"
$ws.Range("F161").Value = "https://www.youtube.com/watch?v=Kn3OJ2njI6Y"
$ws.Range("B162").Value = "PAGANI"
$ws.Range("C162").Value = "BLANKO"
$ws.Range("D162").Value = "PAGANI"
$ws.Range("E162").Value = "'2024-05-03"
$ws.Range("F162").Value = "https://www.youtube.com/watch?v=mtDMUPmmsnk"
$ws.Range("B163").Value = "Alto Troyanaje"
$ws.Range("C163").Value = "Waldokinc El Troyano"
$ws.Range("D163").Value = "Nivel 2"
$ws.Range("E163").Value = "'2013-04-30"
$ws.Range("F163").Value = "https://www.youtube.com/watch?v=bjaa8L37xbM"
$ws.Range("B164").Value = "Ve y Diles"
$ws.Range("C164").Value = "Alex Ponce, Sebastian Llosa"
$ws.Range("D164").Value = "Ve y Diles"
$ws.Range("E164").Value = "'2024-05-09"
$ws.Range("F164").Value = "https://www.youtube.com/watch?v=Qhw83uv0RUs"
$ws.Range("B165").Value = "VIP"
$ws.Range("C165").Value = "Milo Mae"
$ws.Range("D165").Value = "LO TOP, LO EXCLUSIVE"
$ws.Range("E165").Value = "'2024-02-16"
$ws.Range("F165").Value = "https://www.youtube.com/watch?v=_vnr3ZtQRBg"
$ws.Range("B166").Value = "Una y otra vez"
$ws.Range("C166").Value = "Vincez, Milo Mae, Deep Nao"
$ws.Range("D166").Value = "Una y otra vez"
$ws.Range("E166").Value = "'2024-05-03"
$ws.Range("F166").Value = "https://www.youtube.com/watch?v=bQWS4WGiLhA"
$ws.Range("B167").Value = "PLAYER"
$ws.Range("C167").Value = "Carlos Corté$"
$ws.Range("D167").Value = "PLAYER"
$ws.Range("E167").Value = "'2023-09-01"
$ws.Range("F167").Value = "https://www.youtube.com/watch?v=X3qfOg9X7hQ"
$ws.Range("B168").Value = "CONTIGO"
$ws.Range("D168").Value = "PACK VERANO 002"
$ws.Range("E168").Value = "'2022-08-26"
$ws.Range("F168").Value = "https://www.youtube.com/watch?v=TUspIUD8t_s"
$ws.Range("B169").Value = "SEXATON"
$ws.Range("C169").Value = "Milo Mae, Carlos Corté$, Kenny Die, Vyse"
$ws.Range("D169").Value = "LO TOP, LO EXCLUSIVE"
$ws.Range("E169").Value = "'2024-02-16"
$ws.Range("F169").Value = "https://www.youtube.com/watch?v=yr_vgcdAyKk"
$ws.Range("B170").Value = "Tu para Mi"
$ws.Range("C170").Value = "Kenny Die"
$ws.Range("D170").Value = "Tu para Mi"
$ws.Range("E170").Value = "'2022-06-16"
$ws.Range("F170").Value = "https://www.youtube.com/watch?v=P59nEZ-8mDo"
$ws.Range("B171").Value = "Titerito"
$ws.Range("C171").Value = "Ecby, Milo Mae, Strauss"
$ws.Range("D171").Value = "Titerito"
$ws.Range("E171").Value = "'2024-06-13"
$ws.Range("F171").Value = "https://www.youtube.com/watch?v=I6THGFRfBDg"
$ws.Range("B172").Value = "Casa Blanca"
$ws.Range("C172").Value = "Mersa"
$ws.Range("D172").Value = "Casa Blanca"
$ws.Range("E172").Value = "'2024-02-13"
$ws.Range("F172").Value = "https://www.youtube.com/watch?v=4eBG5zOrJW4"
$ws.Range("B173").Value = "Lunares"
$ws.Range("C173").Value = "Mersa, Mr Reo, Kiff"
$ws.Range("D173").Value = "Lunares"
$ws.Range("E173").Value = "'2024-05-31"
$ws.Range("F173").Value = "https://www.youtube.com/watch?v=ddbluhWj-EA"
$ws.Range("B174").Value = "Como tu no hay otra"
$ws.Range("C174").Value = "Dayan"
$ws.Range("D174").Value = "Buena Suerte y Adios"
$ws.Range("E174").Value = "'2024-03-29"
$ws.Range("F174").Value = "https://www.youtube.com/watch?v=WeOPkKxP_E0"
$ws.Range("B175").Value = "Noche De Amanecia"
$ws.Range("C175").Value = "Renn"
$ws.Range("D175").Value = "Noche De Amanecia"
$ws.Range("E175").Value = "'2024-07-19"
$ws.Range("F175").Value = "https://www.youtube.com/watch?v=oUXI0Gf0x_w"
